# Uploaded data for participant-2
# Appends the 5 "P2AppTraces/getevent-N" execution rows to the first
# three sheets (GNUCASH-2.1.3, Antennapod-1.6.2.3, ATimeTracker-0.20),
# mirroring the rows that already exist on the fourth sheet
# (GNUCASH-1.0.3).

$wb = $excel.ActiveWorkbook

$values = @(
    "P2AppTraces/getevent-1",
    "P2AppTraces/getevent-2",
    "P2AppTraces/getevent-3",
    "P2AppTraces/getevent-4",
    "P2AppTraces/getevent-5"
)

# ---------------------------------------------------------------------
# Sheet 2: Antennapod-1.6.2.3  -> new rows 9-13 (Execution number 8-12)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate() | Out-Null

$startRow = 9
$startExec = 8
for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $ws2.Cells.Item($row, 1).Value = $startExec + $i
    $ws2.Cells.Item($row, 2).Value = $values[$i]

    # reproduce the existing cell formatting (style index 4) used by
    # the rows directly above this block
    $ws2.Range("B4").Copy() | Out-Null
    $ws2.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null
}
$ws2.Range("A9:B13").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 3: ATimeTracker-0.20  -> new rows 9-13 (Execution number 8-12)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate() | Out-Null

for ($i = 0; $i -lt 5; $i++) {
    $row = $startRow + $i
    $ws3.Cells.Item($row, 1).Value = $startExec + $i
    $ws3.Cells.Item($row, 2).Value = $values[$i]

    $ws3.Range("B4").Copy() | Out-Null
    $ws3.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null
}
$ws3.Range("D29").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet 1: GNUCASH-2.1.3  -> new rows 17-21 (Execution number 16-20)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate() | Out-Null

$row1Start = 17
$exec1Start = 16
for ($i = 0; $i -lt 5; $i++) {
    $row = $row1Start + $i
    $ws1.Cells.Item($row, 1).Value = $exec1Start + $i
    $ws1.Cells.Item($row, 2).Value = $values[$i]

    # column A on this sheet uses style index 2
    $ws1.Range("A12").Copy() | Out-Null
    $ws1.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null

    # column B alternates between style 2 and style 4, matching the
    # same pattern used by rows 12-16 directly above this block
    $ws1.Range("B12").Copy() | Out-Null
    $ws1.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null

    $ws1.Rows.Item($row).RowHeight = 15.75
}

# row 18 (second of the block) uses style 4 on column B, same as row 13
$ws1.Range("B13").Copy() | Out-Null
$ws1.Cells.Item(18, 2).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(18, 2).Value = $values[1]

$ws1.Range("C31").Select() | Out-Null
